# Update profit files after running on 2025-09-04
# Appends a new data row (row 18) to Sheet1: Date + Profit for 09/04/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Range("A18")

# Force text formatting first so Excel doesn't auto-convert the
# "MM/DD/YYYY"-looking string into a date serial number (the sheet's other
# date cells are stored as plain text, not real dates). Resetting the style
# back to "Normal" afterwards drops the temporary NumberFormat so the cell
# keeps the sheet's default (unstyled) look, matching the surrounding rows.
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/04/2025"
$dateCell.Style = "Normal"

$ws.Range("B18").Value = 13529.73

Write-Output "Added row 18 (09/04/2025, 13529.73)"
